# Pilot Gantt (Sprints) - sprint-view rework with hardcoded Gantt bars.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-color the four shared fills (header / engineer / phase / gantt-bar).
#    Colors are OLE (BGR) integers decoded from the target RGB hex values:
#      header blue : 00002060 -> 6299648
#      engineer red: 00C00000 -> 192
#      phase green : 0000B050 -> 5287936
#      bar gray    : 00AEAAAA -> 11184814
#    Setting PatternColor then Color makes fgColor/bgColor match, same as
#    the solid fills already in the workbook.
# ---------------------------------------------------------------------------
function Set-Fill($rangeAddr, $color) {
    $r = $ws.Range($rangeAddr)
    $r.Interior.PatternColor = $color
    $r.Interior.Color = $color
}

Set-Fill "A1" 6299648
Set-Fill "B1:M1" 6299648
Set-Fill "A2:M2" 192
Set-Fill "A16:M16" 192
Set-Fill "A3:M3" 5287936
Set-Fill "A7:M7" 5287936
Set-Fill "A11:M11" 5287936
Set-Fill "A17:M17" 5287936
Set-Fill "A21:M21" 5287936
Set-Fill "A25:M25" 5287936

# ---------------------------------------------------------------------------
# 2) Column widths: A grows to 105, B:M shrink to 12.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 104.08333333333333
$ws.Range("B1:M1").EntireColumn.ColumnWidth = 11.083333333333334

# ---------------------------------------------------------------------------
# 3) Sprint header labels (row 1) - now wrap onto two lines.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Sprint 1`n(W1-W2)"
$ws.Range("C1").Value = "Sprint 2`n(W3-W4)"
$ws.Range("D1").Value = "Sprint 3`n(W5-W6)"
$ws.Range("E1").Value = "Sprint 4`n(W7-W8)"
$ws.Range("F1").Value = "Sprint 5`n(W9-W10)"
$ws.Range("G1").Value = "Sprint 6`n(W11-W12)"
$ws.Range("H1").Value = "Sprint 7`n(W13-W14)"
$ws.Range("I1").Value = "Sprint 8`n(W15-W16)"
$ws.Range("J1").Value = "Sprint 9`n(W17-W18)"
$ws.Range("K1").Value = "Sprint 10`n(W19-W20)"
$ws.Range("L1").Value = "Sprint 11`n(W21-W22)"
$ws.Range("M1").Value = "Sprint 12`n(W23-W24)"

# ---------------------------------------------------------------------------
# 4) Phase header labels get "Phase N: ... (Est. Months X-Y)" wording.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Phase 1: Discovery, Analysis & Planning (Est. Months 1-2)"
$ws.Range("A7").Value = "Phase 2: Migration, Automation Development & Initial Integration (Est. Months 3-4)"
$ws.Range("A11").Value = "Phase 3: Refinement, Reporting & Knowledge Transfer Preparation (Est. Months 5-6)"
$ws.Range("A17").Value = "Phase 1: Assessment, Strategy Definition & Foundational Setup (Est. Months 1-2)"
$ws.Range("A21").Value = "Phase 2: Implementation, Coaching & CI/CD Integration (Est. Months 3-4)"
$ws.Range("A25").Value = "Phase 3: Optimization, Standardization & Knowledge Dissemination (Est. Months 5-6)"

# ---------------------------------------------------------------------------
# 5) Task rows: append explicit week ranges to the label and paint the
#    matching Gantt bar cells with the gray fill (style carried by s=6).
# ---------------------------------------------------------------------------
$bar = 11184814

$ws.Range("A4").Value = "1. Deep Dive into Existing UAT Processes & Test Assets (W1-W3)"
Set-Fill "B4:C4" $bar

$ws.Range("A5").Value = "2. Identify & Prioritize UAT Scenarios for Automation (W3-W6)"
Set-Fill "C5:D5" $bar

$ws.Range("A6").Value = "3. Master BDD Tooling & Methodology (W2-W5)"
Set-Fill "B6:D6" $bar

$ws.Range("A8").Value = "4. Convert Selected UAT Scenarios to BDD (Gherkin) (W7-W10)"
Set-Fill "E8:F8" $bar

$ws.Range("A9").Value = "5. Develop Automated Test Scripts using Playwright (W9-W16)"
Set-Fill "F9:I9" $bar

$ws.Range("A10").Value = "6. Setup & Test Execution in DT2 Environment (W15-W18)"
Set-Fill "I10:J10" $bar

$ws.Range("A12").Value = "7. Iterate and Refine Automated UAT Suite (W17-W24)"
Set-Fill "J12:M12" $bar

$ws.Range("A13").Value = "8. Establish Automated UAT Reporting (W19-W22)"
Set-Fill "K13:L13" $bar

$ws.Range("A14").Value = "9. Document Best Practices & Create Migration Playbook (W20-W24)"
Set-Fill "K14:M14" $bar

$ws.Range("A15").Value = "10. Prepare for Knowledge Sharing & Team Onboarding (W22-W24)"
Set-Fill "L15:M15" $bar

$ws.Range("A18").Value = "1. Baseline Current Engineering Practices & CI/CD Maturity (W1-W3)"
Set-Fill "B18:C18" $bar

$ws.Range("A19").Value = "2. Develop & Communicate Pilot Engineering Practices Adoption Strategy (W2-W4)"
Set-Fill "B19:C19" $bar

$ws.Range("A20").Value = "3. Tooling Onboarding & Environment Preparation (W3-W6)"
Set-Fill "C20:D20" $bar

$ws.Range("A22").Value = "4. Drive Adoption of Unit Testing & Developer-Led Testing (W7-W16)"
Set-Fill "E22:I22" $bar

$ws.Range("A23").Value = "5. Integrate Automated Tests into CI/CD Pipelines (GitHub Actions Focus) (W9-W16)"
Set-Fill "F23:I23" $bar

$ws.Range("A24").Value = "6. Establish & Champion Mocking Practices (Mockito/MockFlow) (W10-W15)"
Set-Fill "F24:I24" $bar

$ws.Range("A26").Value = "7. Refine CI/CD Pipelines (GitHub Actions) and Test Execution Efficiency (W17-W24)"
Set-Fill "J26:M26" $bar

$ws.Range("A28").Value = "9. Facilitate Performance Profiling Setup (W20-W23)"
Set-Fill "K28:M28" $bar

$ws.Range("A29").Value = "10. Prepare for Scaling & Knowledge Transfer (W22-W24)"
Set-Fill "L29:M29" $bar

Write-Host "Gantt sprint view applied"
